$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 5
$ws.Range("B3").Value = 0
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 5
$ws.Range("C4").Value = 20
$ws.Range("E4").Value = 54
$ws.Range("F4").Value = 3
$ws.Range("C5").Value = 8
$ws.Range("D5").Value = 56
$ws.Range("F5").Value = 22
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 6
$ws.Range("G6").Value = 0
$ws.Range("F7").Value = 5
